# Update cached market-price figures in the Gungnir_Profits sheets
# (scheduled runner refresh). Overwrites the price/profit columns
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ,
# LeveProfitNQ/HQ) for the rows whose underlying market data changed.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 96
$ws.Range("H96").Value = 1680.5
$ws.Range("I96").Value = 1555.2
$ws.Range("J96").Value = 1805.8
$ws.Range("K96").Value = 4665.6
$ws.Range("L96").Value = 5417.4
$ws.Range("M96").Value = -3292.6
$ws.Range("N96").Value = -8163.4

$ws = $wb.Worksheets.Item("ARM")
# Row 22
$ws.Range("H22").Value = 11672
$ws.Range("J22").Value = 20000
$ws.Range("L22").Value = 20000
$ws.Range("N22").Value = -20598

# Row 32
$ws.Range("H32").Value = 1816.97
$ws.Range("I32").Value = 1798.0737
$ws.Range("J32").Value = 2176
$ws.Range("K32").Value = 1798.0737
$ws.Range("L32").Value = 2176
$ws.Range("M32").Value = -1511.0737
$ws.Range("N32").Value = -2750

# Row 43
$ws.Range("H43").Value = 6751.6665
$ws.Range("J43").Value = 6751.6665
$ws.Range("L43").Value = 6751.6665
$ws.Range("N43").Value = -7377.6665

# Row 62
$ws.Range("H62").Value = 19499.666
$ws.Range("J62").Value = 19499.666
$ws.Range("L62").Value = 19499.666
$ws.Range("N62").Value = -20747.666

# Row 65
$ws.Range("H65").Value = 19499.666
$ws.Range("J65").Value = 19499.666
$ws.Range("L65").Value = 58498.99800000001
$ws.Range("N65").Value = -64738.99800000001

# Row 76
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

# Row 79
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

# Row 105
$ws.Range("H105").Value = 39975
$ws.Range("J105").Value = 39975
$ws.Range("L105").Value = 39975
$ws.Range("N105").Value = -46963

$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 20280
$ws.Range("I82").Value = 4945
$ws.Range("J82").Value = 43282.5
$ws.Range("K82").Value = 4945
$ws.Range("L82").Value = 43282.5
$ws.Range("M82").Value = -4562
$ws.Range("N82").Value = -44048.5

# Row 85
$ws.Range("H85").Value = 20280
$ws.Range("I85").Value = 4945
$ws.Range("J85").Value = 43282.5
$ws.Range("K85").Value = 4945
$ws.Range("L85").Value = 43282.5
$ws.Range("M85").Value = -3619
$ws.Range("N85").Value = -45934.5

$ws = $wb.Worksheets.Item("CRP")
# Row 33
$ws.Range("H33").Value = 16180
$ws.Range("I33").Value = 2206.6667
$ws.Range("J33").Value = 23166.666
$ws.Range("K33").Value = 2206.6667
$ws.Range("L33").Value = 23166.666
$ws.Range("M33").Value = -1827.6667
$ws.Range("N33").Value = -23924.666

# Row 99
$ws.Range("H99").Value = 90921200
$ws.Range("I99").Value = 166685380
$ws.Range("K99").Value = 166685380
$ws.Range("M99").Value = -166683882

# Row 126
$ws.Range("H126").Value = 90921200
$ws.Range("I126").Value = 166685380
$ws.Range("K126").Value = 500056140
$ws.Range("M126").Value = -500053670

# Row 132
$ws.Range("H132").Value = 8773426
$ws.Range("I132").Value = 1277.45
$ws.Range("J132").Value = 18520256
$ws.Range("K132").Value = 3832.35
$ws.Range("L132").Value = 55560768
$ws.Range("M132").Value = -1302.35
$ws.Range("N132").Value = -55565828

$ws = $wb.Worksheets.Item("GSM")
# Row 6
$ws.Range("H6").Value = 15361.8
$ws.Range("J6").Value = 15361.8
$ws.Range("L6").Value = 15361.8
$ws.Range("N6").Value = -15587.8

# Row 16
$ws.Range("H16").Value = 15361.8
$ws.Range("J16").Value = 15361.8
$ws.Range("L16").Value = 15361.8
$ws.Range("N16").Value = -15861.8

# Row 70
$ws.Range("H70").Value = 5197.0312
$ws.Range("I70").Value = 4483.6313
$ws.Range("J70").Value = 6239.6924
$ws.Range("K70").Value = 4483.6313
$ws.Range("L70").Value = 6239.6924
$ws.Range("M70").Value = -4213.6313
$ws.Range("N70").Value = -6779.6924

# Row 73
$ws.Range("H73").Value = 5197.0312
$ws.Range("I73").Value = 4483.6313
$ws.Range("J73").Value = 6239.6924
$ws.Range("K73").Value = 4483.6313
$ws.Range("L73").Value = 6239.6924
$ws.Range("M73").Value = -3547.6313
$ws.Range("N73").Value = -8111.6924

$ws = $wb.Worksheets.Item("LTW")
# Row 62
$ws.Range("H62").Value = 29500
$ws.Range("J62").Value = 29500
$ws.Range("L62").Value = 29500
$ws.Range("N62").Value = -30748

# Row 65
$ws.Range("H65").Value = 29500
$ws.Range("J65").Value = 29500
$ws.Range("L65").Value = 88500
$ws.Range("N65").Value = -94740

# Row 76
$ws.Range("H76").Value = 30000
$ws.Range("J76").Value = 30000
$ws.Range("L76").Value = 30000
$ws.Range("N76").Value = -30676

# Row 79
$ws.Range("H79").Value = 30000
$ws.Range("J79").Value = 30000
$ws.Range("L79").Value = 30000
$ws.Range("N79").Value = -32340

# Row 80
$ws.Range("H80").Value = 39700
$ws.Range("J80").Value = 39700
$ws.Range("L80").Value = 39700
$ws.Range("N80").Value = -41946

# Row 83
$ws.Range("H83").Value = 39700
$ws.Range("J83").Value = 39700
$ws.Range("L83").Value = 119100
$ws.Range("N83").Value = -130332

$ws = $wb.Worksheets.Item("WVR")
# Row 70
$ws.Range("H70").Value = 23744.334
$ws.Range("J70").Value = 23744.334
$ws.Range("L70").Value = 23744.334
$ws.Range("N70").Value = -24374.334

# Row 73
$ws.Range("H73").Value = 23744.334
$ws.Range("J73").Value = 23744.334
$ws.Range("L73").Value = 23744.334
$ws.Range("N73").Value = -25928.334

# Row 82
$ws.Range("H82").Value = 49800
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 49800
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 49800
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -50566

# Row 85
$ws.Range("H85").Value = 49800
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 49800
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 49800
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -52452

# Row 123
$ws.Range("H123").Value = 20083.334
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 20083.334
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 20083.334
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -29883.334

# Row 136
$ws.Range("H136").Value = 1308.9783
$ws.Range("I136").Value = 681.43243
$ws.Range("J136").Value = 3888.889
$ws.Range("K136").Value = 2044.29729
$ws.Range("L136").Value = 11666.667
$ws.Range("M136").Value = 505.70271
$ws.Range("N136").Value = -16766.667

Write-Output "Gungnir_Profits sheets updated"
